# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) holds the same arrears-period value for
# every worker row on the sheet (shared text "2507"). This update bumps
# the period forward to "2508" for all of them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"
